# TC13_Canine_Filter_Breed-Bulldog.xlsx - "Fixed variables and query errors in
# Bread from TC01 to TC30"
#
# The CasesTab query (cell B2 on the "startup" sheet) had an extra
# OPTIONAL MATCH on (co:cohort) plus a trailing `Cohort` output column that
# don't belong in this query (that logic/column lives elsewhere) - drop it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$fixedCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Bulldog']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $fixedCasesQuery

# Row heights settle to the new wrap-height of each query cell once the text
# changed (rows re-measured for the new cell content).
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 244.8

# Selection moved off the (now shorter) CasesTab row back up to B2.
$ws.Activate()
$ws.Range("B2").Select()
